# "Running all the test cases" — on the "Test Cases" sheet (the active sheet),
# flip the Runmode column (D) to "Y" for every test case so that all of them
# get executed, and record the outcome of the first two runs in the Results
# column (E): TestCase_B1 passed, TestCase_B2 failed (the rest stay "SKIP").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Runmode column: set every data row (2-25) to "Y" (run).
$ws.Range("D2:D25").Value = "Y"

# Results for the two test cases that were actually executed.
$ws.Range("E2").Value = "PASS"
$ws.Range("E3").Value = "FAIL"

# Leave the selection on the Runmode column, starting at D2.
$ws.Range("D2:D25").Select()
